# Daily data updated on March 23.
$wb = $excel.ActiveWorkbook

# ---------- longform sheet: add YK/NT columns + new date row ----------
$long = $wb.Worksheets.Item("longform")
$long.Activate()

# New header cells for YK_conf, YK_prob, YK_deaths, NT_conf, NT_prob, NT_deaths
$headers = @("YK_conf","YK_prob","YK_deaths","NT_conf","NT_prob","NT_deaths")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $long.Cells.Item(1, 37 + $i).Value = $headers[$i]
}

# New row 30: data_collection / source / date
$long.Cells.Item(30, 1).Value = "live"
$long.Cells.Item(30, 2).Value = "govt_canada_ph"
$long.Cells.Item(30, 3).Value = 43913.458333333336
$long.Cells.Item(30, 3).NumberFormat = "m/d/yy h:mm"

$row30 = @(424,0,10,259,0,1,33,19,0,11,9,0,425,0,5,221,0,4,3,6,0,9,8,0,28,0,0,3,0,0,13,0,0,2,0,0,1,0,0)
for ($i = 0; $i -lt $row30.Length; $i++) {
    $long.Cells.Item(30, 4 + $i).Value = $row30[$i]
}

$long.Range("C30").Select()

# ---------- shortform sheet: add broken-out YK/NT rows for the new date ----------
$short = $wb.Worksheets.Item("shortform")
$short.Activate()

$shortData = @(
    ,("BC","conf",424)
    ,("BC","prob",0)
    ,("BC","deaths",10)
    ,("AB","conf",259)
    ,("AB","prob",0)
    ,("AB","deaths",1)
    ,("SK","conf",33)
    ,("SK","prob",19)
    ,("SK","deaths",0)
    ,("MB","conf",11)
    ,("MB","prob",9)
    ,("MB","deaths",0)
    ,("ON","conf",425)
    ,("ON","prob",0)
    ,("ON","deaths",5)
    ,("QC","conf",221)
    ,("QC","prob",0)
    ,("QC","deaths",4)
    ,("NL","conf",3)
    ,("NL","prob",6)
    ,("NL","deaths",0)
    ,("NB","conf",9)
    ,("NB","prob",8)
    ,("NB","deaths",0)
    ,("NS","conf",28)
    ,("NS","prob",0)
    ,("NS","deaths",0)
    ,("PEI","conf",3)
    ,("PEI","prob",0)
    ,("PEI","deaths",0)
    ,("Repat","conf",13)
    ,("Repat","prob",0)
    ,("Repat","deaths",0)
    ,("YK","conf",2)
    ,("YK","prob",0)
    ,("YK","deaths",0)
    ,("NT","conf",1)
    ,("NT","prob",0)
    ,("NT","deaths",0)
)

$r = 249
foreach ($item in $shortData) {
    $prov = $item[0]
    $ctype = $item[1]
    $val = $item[2]
    $short.Cells.Item($r, 1).Value = "live"
    $short.Cells.Item($r, 2).Value = "govt_canada_ph"
    $short.Cells.Item($r, 3).Value = 43913.458333333336
    $short.Cells.Item($r, 3).NumberFormat = "m/d/yy h:mm"
    $short.Cells.Item($r, 4).Value = $prov
    $short.Cells.Item($r, 5).Value = $ctype
    $short.Cells.Item($r, 6).Value = $val
    $r++
}

$short.Range("C249:C287").Select()

